$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The data is held in an Excel Table (ListObject) named "Table1" spanning
# D4:J63. Add a new row at the end of the table for the latest report entry.
$tbl = $ws.ListObjects.Item("Table1")
$lastRow = $tbl.ListRows.Item($tbl.ListRows.Count).Range
$newRow = $tbl.ListRows.Add()

$dataRange = $newRow.Range

# Match formatting of the previous last row (Excel normally extends the
# table's row formatting automatically when a new row is appended).
$lastRow.Copy()
$dataRange.PasteSpecial(-4122)
$dataRange.RowHeight = $lastRow.RowHeight

# D: Fecha (stored as text, matching the existing rows' text-formatted dates)
$dataRange.Item(1).Value = "23/7/2027"

# E: Imagenes sin etiquetar
$dataRange.Item(2).Value = 380

# F: Imagenes etiquetadas sin revisar
$dataRange.Item(3).Value = 950

# G: Imagenes rechazadas
$dataRange.Item(4).Value = 0

# H: Imagenes etiquetadas y revisadas, faltando de subir
$dataRange.Item(5).Value = 0

# I: Imagenes etiquetadas, revisadas y subidas
$dataRange.Item(6).Value = 1012

# J: Notas
$dataRange.Item(7).Value = "N/A"

# Update the active selection as recorded after the edit
$ws.Range("E65").Select()
